$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row updates
$ws.Range("B1").Value = "sparsity_necessary"
$ws.Range("C1").Value = "necessary explanation rate"

# CoDy section (rows 2-13): selection strategy "1-best" -> "1-delta"
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 5).Value = "1-delta"
}

# Greedy section (rows 23-35): selection strategy "recent" -> "temporal"
for ($r = 23; $r -le 35; $r++) {
    $ws.Cells.Item($r, 5).Value = "temporal"
}
